$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Lista Clientes")
$ws2 = $wb.Worksheets.Item("Catalogo códigos")

# --- Header row: document-number & ubigeo columns become Text-formatted ---
$ws1.Range("C1").NumberFormat = "@"

# --- Data rows: document numbers / ubigeo normalized with Text format ---
$ws1.Range("C2").NumberFormat = "@"

$ws1.Range("C3").NumberFormat = "@"
$ws1.Range("C3").Value = "0235265"

$ws1.Range("V3").NumberFormat = "@"
$ws1.Range("V3").Value = "023652"

$ws1.Range("V2").NumberFormat = "@"
$ws1.Range("V2").Value = "125632"

$ws1.Range("V1").NumberFormat = "@"

# --- Active sheet / selection ---
$ws1.Activate()
$ws1.Range("C3").Select()
